$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap data between row 2 and row 6 for columns D, M, N, O, P, R, S
$d2 = $ws.Range("D2").Value2
$m2 = $ws.Range("M2").Value2
$n2 = $ws.Range("N2").Value2
$o2 = $ws.Range("O2").Value2
$p2 = $ws.Range("P2").Value2
$r2 = $ws.Range("R2").Value2
$s2 = $ws.Range("S2").Value2

$d6 = $ws.Range("D6").Value2
$m6 = $ws.Range("M6").Value2
$n6 = $ws.Range("N6").Value2
$o6 = $ws.Range("O6").Value2
$p6 = $ws.Range("P6").Value2
$r6 = $ws.Range("R6").Value2
$s6 = $ws.Range("S6").Value2

$ws.Range("D2").Value2 = $d6
$ws.Range("M2").Value2 = $m6
$ws.Range("N2").Value2 = $n6
$ws.Range("O2").Value2 = $o6
$ws.Range("P2").Value2 = $p6
$ws.Range("R2").Value2 = $r6
$ws.Range("S2").Value2 = $s6

$ws.Range("D6").Value2 = $d2
$ws.Range("M6").Value2 = $m2
$ws.Range("N6").Value2 = $n2
$ws.Range("O6").Value2 = $o2
$ws.Range("P6").Value2 = $p2
$ws.Range("R6").Value2 = $r2
$ws.Range("S6").Value2 = $s2
